$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Groups") updated values
$ws.Range("C2").Value = 0.44115549721919844
$ws.Range("D2").Value = 0.1470518324063995
$ws.Range("E2").Value = 4.884781506196352
$ws.Range("G2").Value = 0.0041

# Row 3 ("Residuals") updated values
$ws.Range("C3").Value = 4.334986905690751
$ws.Range("D3").Value = 0.030104075733963547
